$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 564
$wsExhibit.Range("F4").Value = 338
$wsExhibit.Range("F5").Value = 403
$wsExhibit.Range("F6").Value = 257
$wsExhibit.Range("F7").Value = 2375
$wsExhibit.Range("F8").Value = 400
$wsExhibit.Range("F9").Value = 6073
$wsExhibit.Range("F10").Value = 155
$wsExhibit.Range("F11").Value = 390

# Sheet "全部类型" (all types) - update "想去人数" (F column) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 564
$wsAll.Range("F4").Value = 338
$wsAll.Range("F5").Value = 403
$wsAll.Range("F6").Value = 257
$wsAll.Range("F9").Value = 2375
$wsAll.Range("F10").Value = 400
$wsAll.Range("F11").Value = 6073
$wsAll.Range("F12").Value = 155
$wsAll.Range("F13").Value = 390
